$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.188888072967529
$ws.Range("B1").Value = 2.362487316131592
$ws.Range("C1").Value = 4.224532127380371
$ws.Range("D1").Value = 2.879838705062866
$ws.Range("E1").Value = 1.12792444229126
